# "more overhaul of spreadsheet"
#
# Changes:
#  - J2: rewrite formula to use absolute refs for C2/B2/G2 and add the
#        launcher-X offset ($D2)
#  - J3:J12: new formula (absolute $C$2, relative row refs) - one shared
#        fill so the whole block computes together
#  - M1: new header "Xf" (same text as J1, reuses the shared string)
#  - M2: new formula, same shape as J2 but using $C$3 (second V0) instead
#        of $C$2
#  - M3:M11: new formula filled across the same span used by H/I/K/L
#  - Row 11 (H11,I11,J11,K11,L11,M11): the old "extra" data row is wiped
#        back to blank cells (formatting only, no value/formula)
#  - Active selection moves to I10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J column (Xf computed from V0 in C2) ---------------------------------
$ws.Range("J2").Formula = '= 2*$C$2^2*COS($B2)*SIN($B2)/$G$2+$D$2'

# Fill J3:J12 in one shot so the engine treats it as a single formula block.
$ws.Range("J3:J12").Formula = '= 2*$C$2^2*COS($B3)*SIN($B3)/$G$2+$D$2'

# Match the number formatting (bold, 2-decimal) used by the rest of column J.
$ws.Range("J3:J12").NumberFormat = "0.00"
$ws.Range("J3:J12").Font.Bold = $true

# --- M column (new - Xf computed from V0 in C3) ---------------------------
$ws.Range("M1").Value = "Xf"

$ws.Range("M2").Formula = '= 2*$C$3^2*COS($B2)*SIN($B2)/$G$2+$D$2'
$ws.Range("M2").NumberFormat = "0.00"
$ws.Range("M2").Font.Bold = $true

$ws.Range("M3:M11").Formula = '= 2*$C$3^2*COS($B3)*SIN($B3)/$G$2+$D$2'
$ws.Range("M3:M11").NumberFormat = "0.00"
$ws.Range("M3:M11").Font.Bold = $true

# --- Wipe the old extra data row (row 11) back to blank -------------------
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("J12").ClearContents()

# --- Selection -------------------------------------------------------------
$ws.Range("I10").Select() | Out-Null
